$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.957.37'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '1.893.60'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").Value = "'0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = "'0.7748"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("D6").Value = "'243.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").Value = "'0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = "'0.3127"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").Value = "'25.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.01%  '

$ws.Range("D10").Value = "'0.07254"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.62%  '

$ws.Range("D11").Value = "'0.08708"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.77%  '

$ws.Range("D12").Value = '1.983.71'
$ws.Range("E12").Value = '  +11.56%  '

$ws.Range("D13").Value = "'0.7722"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.53%  '

$ws.Range("D14").Value = "'5.416"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.82%  '

$ws.Range("D15").Value = "'94.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.57%  '

$ws.Range("D16").Value = "'6.216"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.19%  '

$ws.Range("D17").Value = '30.243.94'
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").Value = '2.328.22'
$ws.Range("E19").Value = '  +8.62%  '

$ws.Range("D20").Value = "'245.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("D21").Value = "'0.000007892"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.69%  '

$ws.Range("D22").Value = "'8.169"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").Value = "'0.9996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("D25").Value = "'0.1595"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.16%  '

$ws.Range("D26").Value = "'9.533"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.63%  '

$ws.Range("D27").Value = "'162.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").Value = "'18.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.88%  '

$ws.Range("D29").Value = "'2.049"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("D30").Value = "'1.431"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.26%  '

$ws.Range("D31").Value = "'1.545"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").Value = "'4.523"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.20%  '

$ws.Range("E33").Value = '  +0.58%  '

$ws.Range("D34").Value = "'0.05442"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.55%  '

$ws.Range("D35").Value = "'1.251"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.95%  '

$ws.Range("D36").Value = "'0.7520"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.38%  '

$ws.Range("D37").Value = "'1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("D38").Value = "'2.686"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.63%  '

$ws.Range("D39").Value = "'0.01984"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.65%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("D41").Value = "'0.4516"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.51%  '

$ws.Range("D42").Value = "'73.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.27%  '

$ws.Range("D43").Value = "'6.087"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.24%  '

$ws.Range("D44").Value = '1.095.68'
$ws.Range("E44").Value = '  -4.22%  '

$ws.Range("D45").Value = '2.239.48'
$ws.Range("E45").Value = '  +7.00%  '

$ws.Range("D46").Value = "'0.8539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.48%  '

$ws.Range("D47").Value = "'0.9996"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.13%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'1.889"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.18%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'103.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").Value = "'7.624"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.70%  '

$ws.Range("D51").Value = "'9.857"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.95%  '
